$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Update the report date on the cover page.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Date: 5 August 2024", $true, $false, $false, $false,
                         $false, $true, 1, $false, "Date: 7 August 2024", 2)

# ---------------------------------------------------------------------------
# 2) Rewrite the "1.5 Combustion Appliance Safety" paragraph: drop the
#    "DELETE UNWANTED" placeholder run (and the stray trailing clause), then
#    add a new paragraph about gas leak detection tests right after it.
# ---------------------------------------------------------------------------
$combustionPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*DELETE UNWANTED*") {
        $combustionPara = $p
        break
    }
}

$combustionRange = $d.Range($combustionPara.Range.Start, $combustionPara.Range.End - 1)
$combustionRange.Text = "We assessed combustion appliances that burn fossil fuels such as propane, heating oil, or kerosene. These include furnaces, boilers, water heaters, and gas ovens. We visually inspected the combustion appliance(s) in your home, as well as conducted combustion safety tests. This included measuring for carbon monoxide and testing that flue gases are properly exhausting from the home."

$combustionRange.InsertParagraphAfter()
$gasLeakPara = $combustionPara.Next()
$gasLeakPara.Style = "BodyText"
$gasLeakPara.Range.Text = "We also performed gas leak detection tests on your propane appliance(s)."

# ---------------------------------------------------------------------------
# 3) Append a brand-new "5 Hey look, a heading!" section (with its own
#    bookmark) at the very end of the document, after "4 Recommendations".
# ---------------------------------------------------------------------------
$headingPara = $d.Paragraphs.Add()
$headingPara.Range.Text = "5 Hey look, a heading!"
$headingPara.Style = "Heading2"

$bodyPara = $d.Paragraphs.Add()
$bodyPara.Range.Text = "lorem ipsum dolor emet…"
$bodyPara.Style = "FirstParagraph"

$newSectionRange = $d.Range($headingPara.Range.Start, $bodyPara.Range.End)
$d.Bookmarks.Add("hey-look-a-heading", $newSectionRange)
